# "more list view improvments"
#
# Adds two new worksheets at the end of the workbook:
#   - "Raw data_discovery - (6)"  (duplicate of "Raw data_discovery - (5)")
#   - "Raw data_outliers - (6)"   (duplicate of "Raw data_outliers - (5)")
#
# Both new sheets are exact content copies of the "- (5)" pair (same header
# row, same data rows/values), matching the next number in the existing
# "Raw data_discovery/outliers - (N)" series.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Raw data_discovery - (6)", copied from "- (5)" -------------
$discoverySrc = $wb.Worksheets.Item("Raw data_discovery - (5)")
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$discoverySrc.Copy([System.Reflection.Missing]::Value, $afterSheet)
$newDiscovery = $wb.Worksheets.Item($wb.Worksheets.Count)
$newDiscovery.Name = "Raw data_discovery - (6)"

# --- Sheet 2: "Raw data_outliers - (6)", copied from "- (5)" --------------
$outliersSrc = $wb.Worksheets.Item("Raw data_outliers - (5)")
$afterSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$outliersSrc.Copy([System.Reflection.Missing]::Value, $afterSheet2)
$newOutliers = $wb.Worksheets.Item($wb.Worksheets.Count)
$newOutliers.Name = "Raw data_outliers - (6)"
